# Regenerate the localization-status report for archive:
#  - flip status text from "Ready for handoff" to "In Translation"
#  - tighten the now-shorter status columns back down to fit the new text

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}

# Column width closest representable to the new autofit width (~13.41 chars)
$newStatusColWidth = 12.5

$wb.Worksheets.Item("Overview").Columns("E:F").ColumnWidth = $newStatusColWidth
$wb.Worksheets.Item("zh-cn").Columns("C:C").ColumnWidth = $newStatusColWidth
$wb.Worksheets.Item("de-de").Columns("C:C").ColumnWidth = $newStatusColWidth
